# 2016 Stronghold Summer Robot Controller Layout.xlsx
# Apply the edits described in the commit:
#   - Remove the Winch and Weave-lift control bindings (that code was
#     pulled out of the project), which clears several cells on both the
#     "X-Box Controller" sheet and the "Joysticks" sheet.
#   - Rename "Fine Drive Control?" to "Fine Drive Control".
#   - Move the button-to-action mapping on the "Joysticks" sheet up now
#     that Winch/Weave rows are gone: Ball Intake/Ball Aim take over the
#     old Winch Down/Winch Up rows.
#   - Update the saved selection/active sheet state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("X-Box Controller")
$ws2 = $wb.Worksheets.Item("Joysticks")

# ---------------------------------------------------------------------
# Sheet 1: "X-Box Controller"
# ---------------------------------------------------------------------

# Winch/Weave lift bindings removed from the code -> clear the cells that
# documented them.
$ws1.Range("B12:C13").Clear()   # Weave Lift Up / Weave Lift Down
$ws1.Range("B17:C18").Clear()   # Winch Down / Winch Up

# "Fine Drive Control?" -> "Fine Drive Control"
$ws1.Range("B22").Value = "Fine Drive Control"
$ws1.Range("C22").Value = "Fine Drive Control"

# ---------------------------------------------------------------------
# Sheet 2: "Joysticks"
# ---------------------------------------------------------------------

# Button 4 used to be Winch Down, now Ball Intake moves up into its place.
$ws2.Range("C11").Value = "Ball Intake"
$ws2.Range("D11").Value = "Ball Intake"

# Button 5 used to be Winch Up, now Ball Aim moves up into its place
# (keeping the left-aligned styling that Ball Aim had before).
$ws2.Range("C12").Value = "Ball Aim"
$ws2.Range("D12").Value = "Ball Aim"
$ws2.Range("C12:D12").HorizontalAlignment = -4131

# Old Ball Intake (button 6) and Ball Aim (button 7) cells are now empty
# (also drops the left-alignment styling that used to live on C14:D14,
# since it moved up to C12:D12 above).
$ws2.Range("C13:D13").Clear()
$ws2.Range("C14:D14").Clear()

# Weave Lift Down (button 10) / Weave Lift Up (button 11) removed.
$ws2.Range("C17:D18").Clear()

# ---------------------------------------------------------------------
# View state: active sheet/tab + saved selections
# ---------------------------------------------------------------------

$ws1.Range("D7").Select()
$ws2.Select()
$ws2.Range("G20").Select()
